$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Structural change: a new "Data (2Byte)" column is inserted before the --
# old column D ("Description" of table 1). This pushes table 1's Description
# column, the gap column, and the whole second ("Host<-FPGA register read")
# table one column to the right (old D..I -> new E..J).
#
# We shift the cell CONTENT by hand (instead of Columns.Insert, which would
# also shift the per-column width metadata) because in the source workbook
# the column widths of the untouched columns did NOT move with the content.

$d1 = $ws.Range("D1").Value2
$f1 = $ws.Range("F1").Value2
$g1 = $ws.Range("G1").Value2
$h1 = $ws.Range("H1").Value2
$i1 = $ws.Range("I1").Value2

$d2 = $ws.Range("D2").Value2
$e2 = $ws.Range("E2").Value2
$f2 = $ws.Range("F2").Value2
$g2 = $ws.Range("G2").Value2
$h2 = $ws.Range("H2").Value2
$i2 = $ws.Range("I2").Value2

$ws.Range("D1:I2").ClearContents()

$ws.Range("E1").Value = $d1
$ws.Range("G1").Value = $f1
$ws.Range("H1").Value = $g1
$ws.Range("I1").Value = $h1
$ws.Range("J1").Value = $i1

$ws.Range("E2").Value = $d2
$ws.Range("F2").Value = $e2
$ws.Range("G2").Value = $f2
$ws.Range("H2").Value = $g2
$ws.Range("I2").Value = $h2
$ws.Range("J2").Value = $i2

# --- New column D: header (row 1) + placeholder (row 2) for table 1 ---------
$ws.Range("D1").Value = "Data (2Byte)"
$ws.Range("D2").Value = "-"

# --- New example rows describing GRAM / IRAM write instructions -------------
$ws.Range("B3").Value = "0x01"
$ws.Range("C3").Value = "Addr"
$ws.Range("D3").Value = "Data"
$ws.Range("E3").Value = "GRAM Write"

$ws.Range("B4").Value = "0x02"
$ws.Range("C4").Value = "Addr"
$ws.Range("D4").Value = "Data"
$ws.Range("E4").Value = "IRAM Write"

# --- Column widths for the columns whose layout actually changed ------------
# (F, G, H keep their original widths untouched - matches source workbook)
$ws.Columns.Item(4).ColumnWidth = 17.83335    # D  -> ~18.61
$ws.Columns.Item(5).ColumnWidth = 25.16665    # E  -> ~25.94
$ws.Columns.Item(9).ColumnWidth = 18.0        # I  -> ~18.76
$ws.Columns.Item(10).ColumnWidth = 25.3333    # J  -> ~26.09

# --- Keep the sheet's used range extending down to row 129, matching the ----
# previously existing (empty, formatting-only) trailing rows of the sheet.
$ws.Range("J129").NumberFormat = "General"

# --- Update the selected cell -------------------------------------------------
$ws.Range("D5").Select()
